$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (row 19) down to the new rows (20-24)
$ws.Range("A19:V19").Copy()
$ws.Range("A20:V24").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Build a 2D array with the new match data (5 rows x 22 columns, A:V)
$newData = New-Object 'object[,]' 5,22

$newData[0,0] = 19
$newData[0,1] = 'gibraltar'
$newData[0,2] = 'national-league'
$newData[0,3] = '2023-2024'
$newData[0,4] = 45226.875
$newData[0,5] = 'St Josephs'
$newData[0,6] = 2
$newData[0,7] = 'Lynx'
$newData[0,8] = 0
$newData[0,9] = 1.45
$newData[0,10] = '27/10/2023 11:47'
$newData[0,11] = 1.46
$newData[0,12] = '27/10/2023 20:59'
$newData[0,13] = 4.34
$newData[0,14] = '27/10/2023 11:47'
$newData[0,15] = 4.38
$newData[0,16] = '27/10/2023 20:59'
$newData[0,17] = 4.9
$newData[0,18] = '27/10/2023 11:47'
$newData[0,19] = 5.18
$newData[0,20] = '27/10/2023 20:59'
$newData[0,21] = 'https://www.betexplorer.com/football/gibraltar/national-league/st-josephs-lynx/bNQHhBwk/'

$newData[1,0] = 20
$newData[1,1] = 'gibraltar'
$newData[1,2] = 'national-league'
$newData[1,3] = '2023-2024'
$newData[1,4] = 45227.6875
$newData[1,5] = 'Europa FC'
$newData[1,6] = 3
$newData[1,7] = 'Glacis United'
$newData[1,8] = 2
$newData[1,9] = 3.01
$newData[1,10] = '28/10/2023 14:17'
$newData[1,11] = 2.39
$newData[1,12] = '28/10/2023 16:28'
$newData[1,13] = 3.69
$newData[1,14] = '28/10/2023 14:17'
$newData[1,15] = 4.6
$newData[1,16] = '28/10/2023 16:28'
$newData[1,17] = 1.96
$newData[1,18] = '28/10/2023 14:17'
$newData[1,19] = 2.13
$newData[1,20] = '28/10/2023 16:28'
$newData[1,21] = 'https://www.betexplorer.com/football/gibraltar/national-league/europa-fc-glacis-united/xEPLiVhe/'

$newData[2,0] = 21
$newData[2,1] = 'gibraltar'
$newData[2,2] = 'national-league'
$newData[2,3] = '2023-2024'
$newData[2,4] = 45227.8125
$newData[2,5] = 'Lions Gibraltar'
$newData[2,6] = 1
$newData[2,7] = 'Europa Point'
$newData[2,8] = 3
$newData[2,9] = 2.4
$newData[2,10] = '28/10/2023 14:32'
$newData[2,11] = 2.49
$newData[2,12] = '28/10/2023 17:59'
$newData[2,13] = 3.42
$newData[2,14] = '28/10/2023 14:32'
$newData[2,15] = 3.83
$newData[2,16] = '28/10/2023 17:59'
$newData[2,17] = 2.43
$newData[2,18] = '28/10/2023 14:32'
$newData[2,19] = 2.26
$newData[2,20] = '28/10/2023 17:59'
$newData[2,21] = 'https://www.betexplorer.com/football/gibraltar/national-league/lions-gibraltar-europa-point/pvPPjk81/'

$newData[3,0] = 22
$newData[3,1] = 'gibraltar'
$newData[3,2] = 'national-league'
$newData[3,3] = '2023-2024'
$newData[3,4] = 45228.6875
$newData[3,5] = 'College 1975 FC'
$newData[3,6] = 2
$newData[3,7] = 'Magpies'
$newData[3,8] = 6
$newData[3,9] = 8.789999999999999
$newData[3,10] = '29/10/2023 13:22'
$newData[3,11] = 13.08
$newData[3,12] = '29/10/2023 16:17'
$newData[3,13] = 6.51
$newData[3,14] = '29/10/2023 13:22'
$newData[3,15] = 8.710000000000001
$newData[3,16] = '29/10/2023 16:17'
$newData[3,17] = 1.19
$newData[3,18] = '29/10/2023 13:22'
$newData[3,19] = 1.1
$newData[3,20] = '29/10/2023 16:15'
$newData[3,21] = 'https://www.betexplorer.com/football/gibraltar/national-league/college-1975-magpies/OUNTk9N7/'

$newData[4,0] = 23
$newData[4,1] = 'gibraltar'
$newData[4,2] = 'national-league'
$newData[4,3] = '2023-2024'
$newData[4,4] = 45228.8125
$newData[4,5] = 'Lincoln Red Imps'
$newData[4,6] = 5
$newData[4,7] = 'Manchester 62'
$newData[4,8] = 0
$newData[4,9] = 1.14
$newData[4,10] = '29/10/2023 13:31'
$newData[4,11] = 1.17
$newData[4,12] = '29/10/2023 17:50'
$newData[4,13] = 7.49
$newData[4,14] = '29/10/2023 13:31'
$newData[4,15] = 7.2
$newData[4,16] = '29/10/2023 17:50'
$newData[4,17] = 9.31
$newData[4,18] = '29/10/2023 13:31'
$newData[4,19] = 8.73
$newData[4,20] = '29/10/2023 17:50'
$newData[4,21] = 'https://www.betexplorer.com/football/gibraltar/national-league/lincoln-red-imps-manchester-62/QZ1CEh08/'

$ws.Range("A20:V24").Value = $newData

Write-Output "Added rows 20-24 (match index 19-23) to Sheet1"
